# Scrum daily standup, day 9
# Update the "Sprint" sheet with today's (day 9, column O) effort entries,
# mark the documentation task as Done, and leave the UI selection/scroll
# state the way the user left it when they finished entering data.

$wb = $excel.ActiveWorkbook

$sprint = $wb.Worksheets.Item("Sprint")
$chart  = $wb.Worksheets.Item("Chart")

# Task in row 7 ("Write documentation for proffesor") is now finished.
$sprint.Range("F7").Value = "Done"

# Day 9 (column O) effort logged against the in-flight tasks.
$sprint.Range("O7").Value = 1
$sprint.Range("O8").Value = 3
$sprint.Range("O12").Value = 2
$sprint.Range("O13").Value = 1

# Leave the cursor where the user last clicked on the Sprint sheet, keeping
# the same scroll position (topLeftCell = B1).
$sprint.Activate()
$sprint.Range("P12").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1

# Scroll/select state on the Chart sheet.
$chart.Activate()
$chart.Range("Q19").Select()
$excel.ActiveWindow.ScrollRow = 2

# Re-activate the Sprint sheet (it was the active tab in the workbook).
$sprint.Activate()
